# Auto-generated Excel COM script to apply Moogle_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1051.9565
$ws.Range("I28").Value = 719.9
$ws.Range("K28").Value = 719.9
$ws.Range("M28").Value = -234.9
$ws.Range("H40").Value = 9835.5
$ws.Range("J40").Value = 11329.125
$ws.Range("L40").Value = 11329.125
$ws.Range("N40").Value = -11679.125
$ws.Range("H80").Value = 1299.069
$ws.Range("J80").Value = 808.2778
$ws.Range("L80").Value = 2424.8334
$ws.Range("N80").Value = -4420.8334
$ws.Range("H83").Value = 1299.069
$ws.Range("J83").Value = 808.2778
$ws.Range("L83").Value = 7274.500199999999
$ws.Range("N83").Value = -17258.5002
$ws.Range("H107").Value = 413.5357
$ws.Range("I107").Value = 376.72726
$ws.Range("J107").Value = 548.5
$ws.Range("K107").Value = 376.72726
$ws.Range("L107").Value = 548.5
$ws.Range("M107").Value = 1543.27274
$ws.Range("N107").Value = -4388.5
$ws.Range("H116").Value = 22090.8
$ws.Range("I116").Value = 37374.43
$ws.Range("J116").Value = 17439.262
$ws.Range("K116").Value = 37374.43
$ws.Range("L116").Value = 17439.262
$ws.Range("M116").Value = -33932.43
$ws.Range("N116").Value = -24323.262
$ws.Range("H137").Value = 1611.16
$ws.Range("I137").Value = 1676.5
$ws.Range("J137").Value = 1349.8
$ws.Range("K137").Value = 5029.5
$ws.Range("L137").Value = 4049.4
$ws.Range("M137").Value = -2479.5
$ws.Range("N137").Value = -9149.4
$ws.Range("H138").Value = 3403.0344
$ws.Range("I138").Value = 3272.8
$ws.Range("J138").Value = 3939.2942
$ws.Range("K138").Value = 9818.400000000001
$ws.Range("L138").Value = 11817.8826
$ws.Range("M138").Value = -4678.400000000001
$ws.Range("N138").Value = -22097.8826
$ws.Range("H141").Value = 1916.0476
$ws.Range("I141").Value = 1744.1052
$ws.Range("J141").Value = 3549.5
$ws.Range("K141").Value = 5232.3156
$ws.Range("L141").Value = 10648.5
$ws.Range("M141").Value = -52.3155999999999
$ws.Range("N141").Value = -21008.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11071.229
$ws.Range("I32").Value = 6980.0186
$ws.Range("K32").Value = 6980.0186
$ws.Range("M32").Value = -6693.0186
$ws.Range("H45").Value = 4169468.2
$ws.Range("I45").Value = 5884470
$ws.Range("J45").Value = 4463.857
$ws.Range("K45").Value = 5884470
$ws.Range("L45").Value = 4463.857
$ws.Range("M45").Value = -5884093
$ws.Range("N45").Value = -5217.857
$ws.Range("H61").Value = 12289.883
$ws.Range("I61").Value = 12841.615
$ws.Range("K61").Value = 12841.615
$ws.Range("M61").Value = -12629.615
$ws.Range("H136").Value = 12289.883
$ws.Range("I136").Value = 12841.615
$ws.Range("K136").Value = 38524.845
$ws.Range("M136").Value = -35974.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4133.1177
$ws.Range("I20").Value = 3296.3333
$ws.Range("K20").Value = 3296.3333
$ws.Range("M20").Value = -3049.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10744.08
$ws.Range("I31").Value = 4306.2354
$ws.Range("J31").Value = 24424.5
$ws.Range("K31").Value = 4306.2354
$ws.Range("L31").Value = 24424.5
$ws.Range("M31").Value = -4011.2354
$ws.Range("N31").Value = -25014.5
$ws.Range("H34").Value = 10744.08
$ws.Range("I34").Value = 4306.2354
$ws.Range("J34").Value = 24424.5
$ws.Range("K34").Value = 4306.2354
$ws.Range("L34").Value = 24424.5
$ws.Range("M34").Value = -4104.2354
$ws.Range("N34").Value = -24828.5
$ws.Range("H58").Value = 4735.7393
$ws.Range("I58").Value = 3382.9375
$ws.Range("J58").Value = 7827.857
$ws.Range("K58").Value = 3382.9375
$ws.Range("L58").Value = 7827.857
$ws.Range("M58").Value = -3179.9375
$ws.Range("N58").Value = -8233.857
$ws.Range("H96").Value = 76149.75
$ws.Range("J96").Value = 76149.75
$ws.Range("L96").Value = 76149.75
$ws.Range("N96").Value = -81641.75
$ws.Range("H99").Value = 4916.077
$ws.Range("I99").Value = 5695.2666
$ws.Range("K99").Value = 5695.2666
$ws.Range("M99").Value = -4197.2666
$ws.Range("H122").Value = 2717.0312
$ws.Range("I122").Value = 1901.2693
$ws.Range("J122").Value = 6252
$ws.Range("K122").Value = 5703.8079
$ws.Range("L122").Value = 18756
$ws.Range("M122").Value = -3253.8079
$ws.Range("N122").Value = -23656
$ws.Range("H126").Value = 4916.077
$ws.Range("I126").Value = 5695.2666
$ws.Range("K126").Value = 17085.7998
$ws.Range("M126").Value = -14615.7998
$ws.Range("H132").Value = 2232.544
$ws.Range("I132").Value = 2004.8302
$ws.Range("K132").Value = 6014.4906
$ws.Range("M132").Value = -3484.4906
$ws.Range("H134").Value = 2702.6843
$ws.Range("I134").Value = 2132.276
$ws.Range("K134").Value = 6396.828
$ws.Range("M134").Value = -3861.828
$ws.Range("H136").Value = 4735.7393
$ws.Range("I136").Value = 3382.9375
$ws.Range("J136").Value = 7827.857
$ws.Range("K136").Value = 10148.8125
$ws.Range("L136").Value = 23483.571
$ws.Range("M136").Value = -7598.8125
$ws.Range("N136").Value = -28583.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 140
$ws.Range("I2").Value = 53.333332
$ws.Range("J2").Value = 226.66667
$ws.Range("K2").Value = 319.999992
$ws.Range("L2").Value = 1360.00002
$ws.Range("M2").Value = -206.999992
$ws.Range("N2").Value = -1586.00002
$ws.Range("H5").Value = 2062.9092
$ws.Range("I5").Value = 1980.25
$ws.Range("J5").Value = 2283.3333
$ws.Range("K5").Value = 5940.75
$ws.Range("L5").Value = 6849.999899999999
$ws.Range("M5").Value = -5828.75
$ws.Range("N5").Value = -7073.999899999999
$ws.Range("H23").Value = 110
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 60
$ws.Range("M23").Value = 175
$ws.Range("H38").Value = 12.416667
$ws.Range("J38").Value = 16
$ws.Range("L38").Value = 48
$ws.Range("N38").Value = -742
$ws.Range("H135").Value = 2062.9092
$ws.Range("I135").Value = 1980.25
$ws.Range("J135").Value = 2283.3333
$ws.Range("K135").Value = 17822.25
$ws.Range("L135").Value = 20549.9997
$ws.Range("M135").Value = -15287.25
$ws.Range("N135").Value = -25619.9997
$ws.Range("H137").Value = 4518.625
$ws.Range("I137").Value = 2808.4
$ws.Range("J137").Value = 7369
$ws.Range("K137").Value = 8425.200000000001
$ws.Range("L137").Value = 22107
$ws.Range("M137").Value = -3325.200000000001
$ws.Range("N137").Value = -32307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2492.4546
$ws.Range("I2").Value = 3364.5715
$ws.Range("J2").Value = 966.25
$ws.Range("K2").Value = 3364.5715
$ws.Range("L2").Value = 966.25
$ws.Range("M2").Value = -3251.5715
$ws.Range("N2").Value = -1192.25
$ws.Range("H113").Value = 3599.8
$ws.Range("I113").Value = 3756.9285
$ws.Range("K113").Value = 3756.9285
$ws.Range("M113").Value = -1586.9285
$ws.Range("H126").Value = 8698.65
$ws.Range("I126").Value = 8614
$ws.Range("K126").Value = 25842
$ws.Range("M126").Value = -23372
$ws.Range("H132").Value = 7214.933
$ws.Range("I132").Value = 4564.636
$ws.Range("J132").Value = 14503.25
$ws.Range("K132").Value = 13693.908
$ws.Range("L132").Value = 43509.75
$ws.Range("M132").Value = -11163.908
$ws.Range("N132").Value = -48569.75
$ws.Range("H136").Value = 49914.95
$ws.Range("J136").Value = 49914.95
$ws.Range("L136").Value = 149744.85
$ws.Range("N136").Value = -154844.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 251100
$ws.Range("I7").Value = 251100
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 251100
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -250988
$ws.Range("H16").Value = 2738.7273
$ws.Range("I16").Value = 2146.5
$ws.Range("J16").Value = 3449.4
$ws.Range("K16").Value = 2146.5
$ws.Range("L16").Value = 3449.4
$ws.Range("M16").Value = -1976.5
$ws.Range("N16").Value = -3789.4
$ws.Range("H40").Value = 4452.143
$ws.Range("J40").Value = 15666
$ws.Range("L40").Value = 15666
$ws.Range("N40").Value = -15938
$ws.Range("H100").Value = 6298.3335
$ws.Range("I100").Value = 5730.6665
$ws.Range("K100").Value = 5730.6665
$ws.Range("M100").Value = -5189.6665
$ws.Range("H126").Value = 251100
$ws.Range("I126").Value = 251100
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 753300
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -750830
$ws.Range("H136").Value = 6498.824
$ws.Range("I136").Value = 4425.4585
$ws.Range("K136").Value = 13276.3755
$ws.Range("M136").Value = -10726.3755
$ws.Range("N7").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 36269
$ws.Range("I51").Value = 29999.5
$ws.Range("J51").Value = 42538.5
$ws.Range("K51").Value = 29999.5
$ws.Range("L51").Value = 42538.5
$ws.Range("M51").Value = -29489.5
$ws.Range("N51").Value = -43558.5
$ws.Range("H52").Value = 2693488
$ws.Range("I52").Value = 4025732
$ws.Range("J52").Value = 28999.75
$ws.Range("K52").Value = 4025732
$ws.Range("L52").Value = 28999.75
$ws.Range("M52").Value = -4025506
$ws.Range("N52").Value = -29451.75
$ws.Range("H98").Value = 60574.75
$ws.Range("J98").Value = 60574.75
$ws.Range("L98").Value = 60574.75
$ws.Range("N98").Value = -66564.75
$ws.Range("H126").Value = 2305.125
$ws.Range("I126").Value = 1205.8572
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 3617.5716
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -1147.5716
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 4898.2925
$ws.Range("I132").Value = 3768.2
$ws.Range("J132").Value = 11304.6
$ws.Range("K132").Value = 11304.6
$ws.Range("M132").Value = -8774.599999999999

Write-Host "Applied changes successfully"
